$d = $word.ActiveDocument

$replacements = @(
    @{ old = "98×65=6370"; new = "32×88=2816" },
    @{ old = "54×25=1350"; new = "52×35=1820" },
    @{ old = "38×86=3268"; new = "97×76=7372" },
    @{ old = "83×22=1826"; new = "31×55=1705" },
    @{ old = "42×16=672";  new = "52×78=4056" },
    @{ old = "48×59=2832"; new = "79×95=7505" },
    @{ old = "36×36=1296"; new = "33×51=1683" },
    @{ old = "11×11=121";  new = "53×52=2756" },
    @{ old = "66×21=1386"; new = "91×47=4277" },
    @{ old = "40×71=2840"; new = "46×99=4554" },
    @{ old = "77×27=2079"; new = "75×86=6450" },
    @{ old = "20×33=660";  new = "95×33=3135" },
    @{ old = "28×16=448";  new = "68×82=5576" },
    @{ old = "63×14=882";  new = "92×66=6072" },
    @{ old = "83×88=7304"; new = "26×37=962" },
    @{ old = "86×30=2580"; new = "79×55=4345" },
    @{ old = "40×12=480";  new = "92×32=2944" },
    @{ old = "49×99=4851"; new = "22×95=2090" },
    @{ old = "83×46=3818"; new = "40×83=3320" },
    @{ old = "80×72=5760"; new = "67×45=3015" },
    @{ old = "15×51=765";  new = "35×43=1505" },
    @{ old = "45×96=4320"; new = "42×17=714" },
    @{ old = "80×91=7280"; new = "26×48=1248" },
    @{ old = "77×81=6237"; new = "47×71=3337" },
    @{ old = "98×96=9408"; new = "83×73=6059" }
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
